# Atualização de bases das ligas, do dia: 25-05-2024 às 15:10
#
# The underlying data source re-ordered a handful of fixtures that share
# the same Date (two matches played on the same matchday were listed in
# the opposite order). For each of the pairs below, the full row of data
# (columns B:AB - everything except the running index in column A) is
# swapped between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(48, 49),
    @(94, 95),
    @(133, 134),
    @(151, 152),
    @(282, 283)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$($rowA):AB$($rowA)")
    $rangeB = $ws.Range("B$($rowB):AB$($rowB)")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}
